# Edit script for Peeksample.pptx slide master placeholder prompt text.
# Wraps each placeholder's text with MARKER_OPENING/MARKER_CLOSING tokens,
# per the target gold OOXML (pseudo-translation + leading/trailing space
# preservation markers). Only ppt/slideMasters/slideMaster1.xml is touched.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# Shape 1 ("Title Placeholder 1"): "Click to edit Master title style"
$titleShape = $m.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Paragraphs(1).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ickclay otay edithay astermay itletay estylay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '

# Shape 2 ("Text Placeholder 2"): 5 paragraphs (lvl 0..4)
#   "Click to edit Master text styles" / "Second level" / "Third level" /
#   "Fourth level" / "Fifth level"
$bodyShape = $m.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyTr.Paragraphs(1).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ickclay otay edithay astermay exttay esstylay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '
$bodyTr.Paragraphs(2).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   econdsay evellay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '
$bodyTr.Paragraphs(3).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   irdthay evellay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '
$bodyTr.Paragraphs(4).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ourthfay evellay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '
$bodyTr.Paragraphs(5).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ifthfay evellay[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '

# Shape 3 ("Date Placeholder 3"): datetimeFigureOut field cached text "3/1/2007"
$dateShape = $m.Shapes.Item(3)
$dateShape.TextFrame.TextRange.Paragraphs(1).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   3/1/2007[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:fld>]   '

# Shape 5 ("Slide Number Placeholder 5"): slidenum field cached text "<#>"
$numShape = $m.Shapes.Item(5)
$numShape.TextFrame.TextRange.Paragraphs(1).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ‹#›[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:fld>]   '
